$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: quarter/period column headers ---
# The oldest reported quarter ("Q2 ended 1399/06") is dropped and all
# quarter columns shift one slot to the left; the newest quarter
# ("Q4 ended 1401/12") is appended in column M.
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9: "date published" labels for each quarter column ---
$ws.Range("D9").Value = "1400-10-30 (2)"
$ws.Range("E9").Value = "1401-02-21 (10)"
$ws.Range("F9").Value = "1401-07-07 (4)"
$ws.Range("G9").Value = "1401-09-07 (5)"
$ws.Range("H9").Value = "1401-10-29 (2)"
$ws.Range("I9").Value = "1402-02-23 (10)"
$ws.Range("J9").Value = "1401-07-07 (2)"
$ws.Range("K9").Value = "1401-09-07 (3)"
$ws.Range("L9").Value = "1401-10-29"
$ws.Range("M9").Value = "1402-02-23 (2)"

# --- Rows 11-27: the income-statement figures themselves shift left by
# one quarter along with the headers, a new quarter of data is appended
# in column M, and a handful of prior-quarter figures are corrected. ---
# Row 11
$ws.Range("D11").Value = 2260008
$ws.Range("E11").Value = 3680202
$ws.Range("F11").Value = 1704944
$ws.Range("G11").Value = 5839350
$ws.Range("H11").Value = 8146415
$ws.Range("I11").Value = 12002785
$ws.Range("J11").Value = 8812065
$ws.Range("K11").Value = 20009824
$ws.Range("L11").Value = 24944661
$ws.Range("M11").Value = 37894630
# Row 12
$ws.Range("D12").Value = -1132437
$ws.Range("E12").Value = -1095691
$ws.Range("F12").Value = -767416
$ws.Range("G12").Value = -2810101
$ws.Range("H12").Value = -3718324
$ws.Range("I12").Value = -5738994
$ws.Range("J12").Value = -5298903
$ws.Range("K12").Value = -12706998
$ws.Range("L12").Value = -16732098
$ws.Range("M12").Value = -20559003
# Row 13
$ws.Range("D13").Value = 1127571
$ws.Range("E13").Value = 2584511
$ws.Range("F13").Value = 937528
$ws.Range("G13").Value = 3029249
$ws.Range("H13").Value = 4428091
$ws.Range("I13").Value = 6263791
$ws.Range("J13").Value = 3513162
$ws.Range("K13").Value = 7302826
$ws.Range("L13").Value = 8212563
$ws.Range("M13").Value = 17335627
# Row 14
$ws.Range("D14").Value = -263075
$ws.Range("E14").Value = -282687
$ws.Range("F14").Value = -286451
$ws.Range("G14").Value = -496201
$ws.Range("H14").Value = -549452
$ws.Range("I14").Value = -756258
$ws.Range("J14").Value = -1055843
$ws.Range("K14").Value = -1698465
$ws.Range("L14").Value = -1721845
$ws.Range("M14").Value = -5371732
# Row 16
$ws.Range("D16").Value = 101731
$ws.Range("E16").Value = 235261
$ws.Range("F16").Value = 15565
$ws.Range("G16").Value = -93524
$ws.Range("H16").Value = 83855
$ws.Range("I16").Value = 211448
$ws.Range("J16").Value = 42415
$ws.Range("K16").Value = 75551
$ws.Range("L16").Value = -475143
$ws.Range("M16").Value = 559096
# Row 17
$ws.Range("D17").Value = 966227
$ws.Range("E17").Value = 2537085
$ws.Range("F17").Value = 666642
$ws.Range("G17").Value = 2439524
$ws.Range("H17").Value = 3962494
$ws.Range("I17").Value = 5718981
$ws.Range("J17").Value = 2499734
$ws.Range("K17").Value = 5679912
$ws.Range("L17").Value = 6015575
$ws.Range("M17").Value = 12522991
# Row 18
$ws.Range("D18").Value = -52414
$ws.Range("E18").Value = -175024
$ws.Range("F18").Value = -249848
$ws.Range("G18").Value = -580703
$ws.Range("H18").Value = -421466
$ws.Range("I18").Value = -701682
$ws.Range("J18").Value = -691689
$ws.Range("K18").Value = -1093372
$ws.Range("L18").Value = -1675504
$ws.Range("M18").Value = -1782092
# Row 19
$ws.Range("D19").Value = -49138
$ws.Range("E19").Value = 32076
$ws.Range("F19").Value = 9222
$ws.Range("G19").Value = 18533
$ws.Range("H19").Value = 247768
$ws.Range("I19").Value = 439818
$ws.Range("J19").Value = 28561
$ws.Range("K19").Value = -248400
$ws.Range("L19").Value = 329346
$ws.Range("M19").Value = 293899
# Row 20
$ws.Range("D20").Value = 864675
$ws.Range("E20").Value = 2394137
$ws.Range("F20").Value = 426016
$ws.Range("G20").Value = 1877354
$ws.Range("H20").Value = 3788796
$ws.Range("I20").Value = 5457117
$ws.Range("J20").Value = 1836606
$ws.Range("K20").Value = 4338140
$ws.Range("L20").Value = 4669417
$ws.Range("M20").Value = 11034798
# Row 21
$ws.Range("D21").Value = -233583
$ws.Range("E21").Value = -62507
$ws.Range("F21").Value = -101178
$ws.Range("G21").Value = -416884
$ws.Range("H21").Value = -577002
$ws.Range("I21").Value = -467328
$ws.Range("J21").Value = -436193
$ws.Range("K21").Value = -720252
$ws.Range("L21").Value = -878594
$ws.Range("M21").Value = -2277108
# Row 22
$ws.Range("D22").Value = 631092
$ws.Range("E22").Value = 2331630
$ws.Range("F22").Value = 324838
$ws.Range("G22").Value = 1460470
$ws.Range("H22").Value = 3211794
$ws.Range("I22").Value = 4989789
$ws.Range("J22").Value = 1400413
$ws.Range("K22").Value = 3617888
$ws.Range("L22").Value = 3790823
$ws.Range("M22").Value = 8757690
# Row 24
$ws.Range("D24").Value = 631092
$ws.Range("E24").Value = 2331630
$ws.Range("F24").Value = 324838
$ws.Range("G24").Value = 1460470
$ws.Range("H24").Value = 3211794
$ws.Range("I24").Value = 4989789
$ws.Range("J24").Value = 1400413
$ws.Range("K24").Value = 3617888
$ws.Range("L24").Value = 3790823
$ws.Range("M24").Value = 8757690
# Row 25
$ws.Range("D25").Value = 631
$ws.Range("E25").Value = 2332
$ws.Range("F25").Value = 325
$ws.Range("G25").Value = 1460
$ws.Range("H25").Value = 3212
$ws.Range("I25").Value = 295
$ws.Range("J25").Value = 175
$ws.Range("K25").Value = 214
$ws.Range("L25").Value = 224
$ws.Range("M25").Value = 518
# Row 26
$ws.Range("I26").Value = 16900000
$ws.Range("K26").Value = 16900000
# Row 27
$ws.Range("D27").Value = 37
$ws.Range("E27").Value = 138
$ws.Range("F27").Value = 19
$ws.Range("G27").Value = 86
$ws.Range("H27").Value = 190
$ws.Range("I27").Value = 295
$ws.Range("J27").Value = 83
$ws.Range("K27").Value = 214
$ws.Range("L27").Value = 224
$ws.Range("M27").Value = 518
